$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Líder del grupo de investigación ..." entries (rows 3 & 9) ---
# Set these first so the new shared string is appended before the
# "Investigador ..." one (matches the order the strings were re-authored).
$ws.Range("E3").Value = "Líder del grupo de investigación CODEC: Ciencias Cognitivas y del Comportamiento (desde 2016)"
$ws.Range("E9").Value = "Líder del grupo de investigación CODEC: Ciencias Cognitivas y del Comportamiento (desde 2016)"

# --- Update the "Investigador ..." entries (rows 2 & 7) ---
$ws.Range("E2").Value = "Investigador en EvoCo: Laboratorio de Evolución y Comportamiento Humano"
$ws.Range("E7").Value = "Investigador en EvoCo: Laboratorio de Evolución y Comportamiento Humano"

# --- Remove the stray fill style that was applied to E15 ---
$ws.Range("E15").Style = "Normal"

# --- Update the window selection / view state ---
$ws.Range("E19").Select()
